# Apply cell value updates to match the refreshed crypto price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.370.97"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "2.661.33"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'602.95"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'178.95"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.174"
$ws.Range("E9").Value = "  +5.70%  "
$ws.Range("D10").Value = "2.660.67"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "'5.05"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "3.147.23"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'0.0000190"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "72.273.98"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "26.59"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "2.664.54"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'12.04"
$ws.Range("E19").Value = "  +4.82%  "
$ws.Range("D20").Value = "'8.01"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'379.82"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'4.20"
$ws.Range("E23").Value = "  +12.83%  "
$ws.Range("D24").Value = "'72.66"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "'4.41"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  +5.08%  "
$ws.Range("D28").Value = "2.807.10"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").Value = "'521.99"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'164.07"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "19.52"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D39").Value = "'1.39"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "0.111"
$ws.Range("E40").Value = "  -6.55%  "
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "'5.10"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.63"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "'39.40"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").Value = "'152.96"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").Value = "'0.549"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").Value = "0.0₆0260"
$ws.Range("E51").Value = "  -2.30%  "
